# Apply the "output generated at c8c62b6" update:
# - update the date heading
# - update each arithmetic problem in the table

$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "2025-09-26 Friday"; New = "2025-09-27 Saturday" },
    @{ Old = "27÷7=";             New = "57÷7=" },
    @{ Old = "75÷5=";             New = "73÷4=" },
    @{ Old = "79÷6=";             New = "28÷3=" },
    @{ Old = "62÷4=";             New = "50÷9=" },
    @{ Old = "23÷9=";             New = "96÷2=" },
    @{ Old = "96÷4=";             New = "77÷9=" },
    @{ Old = "94÷5=";             New = "60÷5=" },
    @{ Old = "95÷2=";             New = "48÷5=" },
    @{ Old = "22÷5=";             New = "60÷2=" },
    @{ Old = "78÷5=";             New = "41÷6=" },
    @{ Old = "92÷8=";             New = "12÷5=" },
    @{ Old = "65÷4=";             New = "63÷2=" },
    @{ Old = "84÷7=";             New = "70÷8=" },
    @{ Old = "24÷7=";             New = "36÷9=" },
    @{ Old = "24÷3=";             New = "22÷2=" },
    @{ Old = "53÷4=";             New = "45÷5=" },
    @{ Old = "50÷3=";             New = "13÷6=" },
    @{ Old = "92÷5=";             New = "97÷9=" },
    @{ Old = "83÷7=";             New = "27÷3=" },
    @{ Old = "81÷8=";             New = "45÷2=" },
    @{ Old = "21÷5=";             New = "86÷2=" },
    @{ Old = "50÷8=";             New = "72÷8=" },
    @{ Old = "49÷5=";             New = "11÷2=" },
    @{ Old = "48÷9=";             New = "34÷2=" },
    @{ Old = "31÷6=";             New = "12÷9=" }
)

foreach ($pair in $replacements) {
    $range = $d.Content
    $range.Find.Execute($pair.Old, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $pair.New, 2)
}

$d.Save()
